# AES Research Paper edits
#
# The only substantive textual changes in the commit are two sentence-level
# wording tweaks (everything else in the recorded diff -- proofErr spell-check
# tags and the run re-segmentation around unchanged paragraphs -- is
# incidental churn Word's proofing pass produces when nearby text is edited;
# the paragraph order/content of the INTRODUCTION...Statement of the Problem
# block is unchanged, it only *looks* like a big move in the line-based diff
# because two sentences inside it changed length and reflowed the
# surrounding lines).
#
# 1. Abstract: "... votes on the canvassing system-level." -> "... votes on
#    the server-level."
# 2. Statement of the Problem question is reworded.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "canvassing system-level.", $true, $false, $false, $false, $false,
    $true, 1, $false, "server-level.", 2) | Out-Null

$d.Content.Find.Execute(
    "How can the Philippine election system prevent electoral fraud specifically for the counting and transmission of the election returns?",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "How can the Philippine automated election system secure the transmission of election returns on the server-level?",
    2) | Out-Null
